$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.282.75'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.897.67'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.73%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.36'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.662'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.97%  '

$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.99'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.81%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.348'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +5.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.92'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +12.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0716'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.176.49'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.07'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.90%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.918.79'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.85%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.694'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.05%  '

$ws.Range("E17").Value = '  +1.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.284.49'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.98'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.72%  '

$ws.Range("E20").Value = '  +2.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '239.82'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.40'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.09%  '

$ws.Range("E23").Value = '  +1.21%  '

$ws.Range("E25").Value = '  +1.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +20.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.11'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.37'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.88%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.27'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.85%  '

$ws.Range("E30").Value = '  +1.50%  '

$ws.Range("E31").Value = '  +2.09%  '

$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.924'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +12.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.06'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.36%  '

$ws.Range("E36").Value = '  -5.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.01'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.02%  '

$ws.Range("E38").Value = '  +1.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.09'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0207'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.02'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0626'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '89.20'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.333.69'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.70'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +36.67%  '

$ws.Range("E47").Value = '  -0.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.76'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.47'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.62'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.086.12'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.79%  '
